$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) First clone the existing last row (126), which currently only has
#    cells A-G and K-AA (no H/I/J/AB/AC), into the two brand-new rows
#    127 and 128. This preserves the per-column cell styles (bold +
#    border on A, date format on E) without creating any new style
#    entries in styles.xml.
# ------------------------------------------------------------------
$ws.Range("A126:G126").Copy($ws.Range("A127:G127")) | Out-Null
$ws.Range("K126:AA126").Copy($ws.Range("K127:AA127")) | Out-Null

$ws.Range("A126:G126").Copy($ws.Range("A128:G128")) | Out-Null
$ws.Range("K126:AA126").Copy($ws.Range("K128:AA128")) | Out-Null

# ------------------------------------------------------------------
# 2) Fill in the new row 127 values (match result still pending, so no
#    FTHG/FTAG/FTR and no PL_AhOver/PL_AhUnder columns).
# ------------------------------------------------------------------
$ws.Range("A127").Value2 = 125
$ws.Range("B127").Value2 = 7011617
$ws.Range("C127").Value2 = "Azerbaijan Premier League"
$ws.Range("D127").Value2 = "Azerbaijan Premier League"
$ws.Range("E127").Value2 = 45361.35416666666
$ws.Range("F127").Value2 = "PFK Turan Tovuz"
$ws.Range("G127").Value2 = "Araz FK"
$ws.Range("K127").Value2 = 2.3
$ws.Range("L127").Value2 = 3.1
$ws.Range("M127").Value2 = 2.8
$ws.Range("N127").Value2 = 2.25
$ws.Range("O127").Value2 = 3.1
$ws.Range("P127").Value2 = 2.875
$ws.Range("Q127").Value2 = -0.25
$ws.Range("R127").Value2 = 2.025
$ws.Range("S127").Value2 = 1.775
$ws.Range("T127").Value2 = 2.25
$ws.Range("U127").Value2 = 1.95
$ws.Range("V127").Value2 = 1.85
$ws.Range("W127").Value2 = 0
$ws.Range("X127").Value2 = 0
$ws.Range("Y127").Value2 = 0
$ws.Range("Z127").Value2 = 0
$ws.Range("AA127").Value2 = 0

# ------------------------------------------------------------------
# 3) Fill in the new row 128 values.
# ------------------------------------------------------------------
$ws.Range("A128").Value2 = 126
$ws.Range("B128").Value2 = 7011618
$ws.Range("C128").Value2 = "Azerbaijan Premier League"
$ws.Range("D128").Value2 = "Azerbaijan Premier League"
$ws.Range("E128").Value2 = 45361.45833333334
$ws.Range("F128").Value2 = "FK Qarabag"
$ws.Range("G128").Value2 = "Zira IK"
$ws.Range("K128").Value2 = 1.4
$ws.Range("L128").Value2 = 4
$ws.Range("M128").Value2 = 7
$ws.Range("N128").Value2 = 1.45
$ws.Range("O128").Value2 = 3.8
$ws.Range("P128").Value2 = 6
$ws.Range("Q128").Value2 = -1
$ws.Range("R128").Value2 = 1.775
$ws.Range("S128").Value2 = 2.025
$ws.Range("T128").Value2 = 2.5
$ws.Range("U128").Value2 = 1.975
$ws.Range("V128").Value2 = 1.825
$ws.Range("W128").Value2 = 0
$ws.Range("X128").Value2 = 0
$ws.Range("Y128").Value2 = 0
$ws.Range("Z128").Value2 = 0
$ws.Range("AA128").Value2 = 0

# ------------------------------------------------------------------
# 4) Now update row 126 itself: the match has since been played, so
#    FTHG/FTAG/FTR (H/I/J) are now populated, the closing odds
#    (N..S) changed, and the Asian-handicap-over/under columns
#    (AB/AC) are now populated too.
# ------------------------------------------------------------------
$ws.Range("H126").Value2 = 0
$ws.Range("I126").Value2 = 1
$ws.Range("J126").Value2 = "A"

$ws.Range("N126").Value2 = 2.05
$ws.Range("O126").Value2 = 3.4
$ws.Range("P126").Value2 = 3
$ws.Range("Q126").Value2 = -0.25
$ws.Range("R126").Value2 = 1.825
$ws.Range("S126").Value2 = 1.975

$ws.Range("U126").Value2 = 1.85
$ws.Range("V126").Value2 = 1.95
$ws.Range("W126").Value2 = -1
$ws.Range("X126").Value2 = -1
$ws.Range("Y126").Value2 = 2
$ws.Range("Z126").Value2 = -1
$ws.Range("AA126").Value2 = 0.9750000000000001
$ws.Range("AB126").Value2 = -1
$ws.Range("AC126").Value2 = 0.95

# ------------------------------------------------------------------
# 5) The worksheet dimension now spans two extra rows.
# ------------------------------------------------------------------
$ws.UsedRange | Out-Null
